# Add a new "API_DATA_FILE" worksheet after the existing "DATA_FILE" sheet
# and populate it with the API test-case data (Page classes for requests).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("DATA_FILE")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "API_DATA_FILE"

# --- data (order chosen so shared-strings get appended in the same
#     sequence as the authored workbook: TC name, JOB header, hyperlink,
#     Developer, Ak)
$ws.Range("A2").Value = "TC_001_POST_SINGLE_USER"
$ws.Range("C1").Value = "JOB"

# Hyperlink first (its TextToDisplay briefly owns the cell text), then
# restore the real cell text - mirrors DATA_FILE's own mismatched
# hyperlink/display-text rows (e.g. D2 / F2).
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:supriya@36", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "supriya@36")
$ws.Range("C2").Value = "Developer"
$ws.Range("B2").Value = "Ak"

# --- remaining header row (reuses existing shared strings from DATA_FILE)
$ws.Range("A1").Value = "TC_NAME"
$ws.Range("B1").Value = "NAME"
$ws.Range("D1").Value = "PRODUCT"
$ws.Range("E1").Value = "EMAIL"
$ws.Range("F1").Value = "NAME"
$ws.Range("G1").Value = "MOBILE"
$ws.Range("H1").Value = "ADDRESS"
$ws.Range("I1").Value = "SEARCH ITEM"
$ws.Range("A1:I1").Font.Bold = $true

# --- hyperlink-style placeholder cells (matches template look copied from
#     the DATA_FILE sheet's EMAIL / NAME columns)
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"

# --- page setup (matches DATA_FILE's printable page settings)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- view state: new sheet becomes the active tab/selection
$ws.Range("D14").Select()
$ws.Activate()

# DATA_FILE's own view scrolls right one column with A1:J5 selected
$dataSheet.Activate()
$dataSheet.Range("A1:J5").Select()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Activate()
